# Apply updated crypto price/volume data to Sheet1.
#
# Every new value is written with a leading apostrophe (a plain PowerShell
# string, not a special escape) so Excel stores it verbatim as TEXT instead
# of auto-coercing digit-and-dot strings like "0.0760" / "2.40" into numbers
# (which would silently drop the significant trailing zero). The apostrophe
# itself is never part of the stored value. Because the quote-prefix trick
# nudges the cell onto a new "quotePrefix" style, .Style is reset back to
# "Normal" right after so the cell keeps its original (unstyled) formatting,
# matching the source workbook where only data rows change, not formatting.
#
# Rows 13/14 (Chainlink <-> WrappedliquidstakedEther2.0) and rows 30/31
# (EthereumClassic <-> Kaspa) are full row-content swaps with refreshed
# price/volume figures, matching the authoritative diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.Value = "'37.215.37"
$cell.Style = "Normal"
$cell = $ws.Range("E2")
$cell.Value = "'  +1.65%  "
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.Value = "'2.058.03"
$cell.Style = "Normal"
$cell = $ws.Range("E3")
$cell.Value = "'  +1.00%  "
$cell.Style = "Normal"

$cell = $ws.Range("E4")
$cell.Value = "'  +0.06%  "
$cell.Style = "Normal"

$cell = $ws.Range("D5")
$cell.Value = "'232.59"
$cell.Style = "Normal"
$cell = $ws.Range("E5")
$cell.Value = "'  +0.42%  "
$cell.Style = "Normal"

$cell = $ws.Range("E6")
$cell.Value = "'  +3.29%  "
$cell.Style = "Normal"

$cell = $ws.Range("E7")
$cell.Value = "'  +0.03%  "
$cell.Style = "Normal"

$cell = $ws.Range("D8")
$cell.Value = "'57.28"
$cell.Style = "Normal"
$cell = $ws.Range("E8")
$cell.Value = "'  +3.54%  "
$cell.Style = "Normal"

$cell = $ws.Range("E9")
$cell.Value = "'  +3.45%  "
$cell.Style = "Normal"

$cell = $ws.Range("D10")
$cell.Value = "'57.85"
$cell.Style = "Normal"
$cell = $ws.Range("E10")
$cell.Value = "'  +1.59%  "
$cell.Style = "Normal"

$cell = $ws.Range("D11")
$cell.Value = "'0.0760"
$cell.Style = "Normal"
$cell = $ws.Range("E11")
$cell.Value = "'  +1.06%  "
$cell.Style = "Normal"

$cell = $ws.Range("E12")
$cell.Value = "'  +1.05%  "
$cell.Style = "Normal"

$cell = $ws.Range("B13")
$cell.Value = "'WrappedliquidstakedEther2.0"
$cell.Style = "Normal"
$cell = $ws.Range("C13")
$cell.Value = "'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$cell.Style = "Normal"
$cell = $ws.Range("D13")
$cell.Value = "'2.359.42"
$cell.Style = "Normal"
$cell = $ws.Range("E13")
$cell.Value = "'  +0.96%  "
$cell.Style = "Normal"

$cell = $ws.Range("B14")
$cell.Value = "'Chainlink"
$cell.Style = "Normal"
$cell = $ws.Range("C14")
$cell.Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$cell.Style = "Normal"
$cell = $ws.Range("D14")
$cell.Value = "'14.68"
$cell.Style = "Normal"
$cell = $ws.Range("E14")
$cell.Value = "'  +2.67%  "
$cell.Style = "Normal"

$cell = $ws.Range("D15")
$cell.Value = "'20.86"
$cell.Style = "Normal"
$cell = $ws.Range("E15")
$cell.Value = "'  +4.14%  "
$cell.Style = "Normal"

$cell = $ws.Range("D16")
$cell.Value = "'0.780"
$cell.Style = "Normal"
$cell = $ws.Range("E16")
$cell.Value = "'  +2.63%  "
$cell.Style = "Normal"

$cell = $ws.Range("D17")
$cell.Value = "'5.16"
$cell.Style = "Normal"
$cell = $ws.Range("E17")
$cell.Value = "'  +0.14%  "
$cell.Style = "Normal"

$cell = $ws.Range("D18")
$cell.Value = "'2.058.23"
$cell.Style = "Normal"
$cell = $ws.Range("E18")
$cell.Value = "'  +1.02%  "
$cell.Style = "Normal"

$cell = $ws.Range("D19")
$cell.Value = "'37.182.62"
$cell.Style = "Normal"
$cell = $ws.Range("E19")
$cell.Value = "'  +1.22%  "
$cell.Style = "Normal"

$cell = $ws.Range("D20")
$cell.Value = "'6.39"
$cell.Style = "Normal"
$cell = $ws.Range("E20")
$cell.Value = "'  +9.18%  "
$cell.Style = "Normal"

$cell = $ws.Range("E21")
$cell.Value = "'  +2.38%  "
$cell.Style = "Normal"

$cell = $ws.Range("E22")
$cell.Value = "'  +1.62%  "
$cell.Style = "Normal"

$cell = $ws.Range("D23")
$cell.Value = "'225.94"
$cell.Style = "Normal"
$cell = $ws.Range("E23")
$cell.Value = "'  +2.37%  "
$cell.Style = "Normal"

$cell = $ws.Range("D24")
$cell.Value = "'0.999"
$cell.Style = "Normal"
$cell = $ws.Range("E24")
$cell.Value = "'  -0.08%  "
$cell.Style = "Normal"

$cell = $ws.Range("D25")
$cell.Value = "'2.40"
$cell.Style = "Normal"
$cell = $ws.Range("E25")
$cell.Value = "'  +0.16%  "
$cell.Style = "Normal"

$cell = $ws.Range("E26")
$cell.Value = "'  +1.16%  "
$cell.Style = "Normal"

$cell = $ws.Range("D27")
$cell.Value = "'165.69"
$cell.Style = "Normal"
$cell = $ws.Range("E27")
$cell.Value = "'  +1.85%  "
$cell.Style = "Normal"

$cell = $ws.Range("E28")
$cell.Value = "'  +7.64%  "
$cell.Style = "Normal"

$cell = $ws.Range("D29")
$cell.Value = "'8.79"
$cell.Style = "Normal"
$cell = $ws.Range("E29")
$cell.Value = "'  +0.72%  "
$cell.Style = "Normal"

$cell = $ws.Range("B30")
$cell.Value = "'Kaspa"
$cell.Style = "Normal"
$cell = $ws.Range("C30")
$cell.Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$cell.Style = "Normal"
$cell = $ws.Range("D30")
$cell.Value = "'0.127"
$cell.Style = "Normal"
$cell = $ws.Range("E30")
$cell.Value = "'  +1.24%  "
$cell.Style = "Normal"

$cell = $ws.Range("B31")
$cell.Value = "'EthereumClassic"
$cell.Style = "Normal"
$cell = $ws.Range("C31")
$cell.Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$cell.Style = "Normal"
$cell = $ws.Range("D31")
$cell.Value = "'19.03"
$cell.Style = "Normal"
$cell = $ws.Range("E31")
$cell.Value = "'  +0.73%  "
$cell.Style = "Normal"

$cell = $ws.Range("D32")
$cell.Value = "'0.117"
$cell.Style = "Normal"
$cell = $ws.Range("E32")
$cell.Value = "'  +0.69%  "
$cell.Style = "Normal"

$cell = $ws.Range("D33")
$cell.Value = "'4.46"
$cell.Style = "Normal"
$cell = $ws.Range("E33")
$cell.Value = "'  +2.46%  "
$cell.Style = "Normal"

$cell = $ws.Range("E34")
$cell.Value = "'  +2.07%  "
$cell.Style = "Normal"

$cell = $ws.Range("D35")
$cell.Value = "'4.60"
$cell.Style = "Normal"
$cell = $ws.Range("E35")
$cell.Value = "'  +7.82%  "
$cell.Style = "Normal"

$cell = $ws.Range("D36")
$cell.Value = "'2.50"
$cell.Style = "Normal"
$cell = $ws.Range("E36")
$cell.Value = "'  +0.84%  "
$cell.Style = "Normal"

$cell = $ws.Range("E37")
$cell.Value = "'  +0.18%  "
$cell.Style = "Normal"

$cell = $ws.Range("D38")
$cell.Value = "'3.26"
$cell.Style = "Normal"
$cell = $ws.Range("E38")
$cell.Value = "'  +1.56%  "
$cell.Style = "Normal"

$cell = $ws.Range("E39")
$cell.Value = "'  -0.87%  "
$cell.Style = "Normal"

$cell = $ws.Range("E41")
$cell.Value = "'  +0.21%  "
$cell.Style = "Normal"

$cell = $ws.Range("D42")
$cell.Value = "'4.44"
$cell.Style = "Normal"
$cell = $ws.Range("E42")
$cell.Value = "'  +0.68%  "
$cell.Style = "Normal"

$cell = $ws.Range("D43")
$cell.Value = "'1.473.20"
$cell.Style = "Normal"
$cell = $ws.Range("E43")
$cell.Value = "'  -0.06%  "
$cell.Style = "Normal"

$cell = $ws.Range("D44")
$cell.Value = "'96.37"
$cell.Style = "Normal"
$cell = $ws.Range("E44")
$cell.Value = "'  +2.73%  "
$cell.Style = "Normal"

$cell = $ws.Range("E45")
$cell.Value = "'  +5.86%  "
$cell.Style = "Normal"

$cell = $ws.Range("D46")
$cell.Value = "'0.0930"
$cell.Style = "Normal"
$cell = $ws.Range("E46")
$cell.Value = "'  -1.37%  "
$cell.Style = "Normal"

$cell = $ws.Range("E47")
$cell.Value = "'  +3.50%  "
$cell.Style = "Normal"

$cell = $ws.Range("E48")
$cell.Value = "'  +1.80%  "
$cell.Style = "Normal"

$cell = $ws.Range("D49")
$cell.Value = "'15.11"
$cell.Style = "Normal"
$cell = $ws.Range("E49")
$cell.Value = "'  -2.93%  "
$cell.Style = "Normal"

$cell = $ws.Range("E50")
$cell.Value = "'  +3.16%  "
$cell.Style = "Normal"

$cell = $ws.Range("E51")
$cell.Value = "'  +1.73%  "
$cell.Style = "Normal"
